$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 104 (id=102): becomes what row 105 used to contain (match swapped with row 105)
$ws.Range("B104").Value = 7127370
$ws.Range("E104").Value = "Macarthur FC"
$ws.Range("F104").Value = "Wellington Phoenix"
$ws.Range("H104").Value = 2
$ws.Range("I104").Value = "A"
$ws.Range("J104").Value = 2.4
$ws.Range("L104").Value = 2.625
$ws.Range("M104").Value = 2.375
$ws.Range("N104").Value = 3.8
$ws.Range("O104").Value = 2.75
$ws.Range("P104").Value = 0
$ws.Range("Q104").Value = 1.8
$ws.Range("R104").Value = 2.05
$ws.Range("S104").Value = 3
$ws.Range("T104").Value = 1.9
$ws.Range("U104").Value = 1.95
$ws.Range("V104").Value = -1
$ws.Range("X104").Value = 1.75
$ws.Range("Y104").Value = -1
$ws.Range("Z104").Value = 1.05
$ws.Range("AA104").Value = 0
$ws.Range("AB104").Value = 0

# Row 105 (id=103): becomes what row 104 used to contain
$ws.Range("B105").Value = 7127374
$ws.Range("E105").Value = "Central Coast Mariners"
$ws.Range("F105").Value = "Western Sydney Wanderers"
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = "H"
$ws.Range("J105").Value = 1.909
$ws.Range("L105").Value = 3.6
$ws.Range("M105").Value = 2.15
$ws.Range("N105").Value = 3.6
$ws.Range("O105").Value = 3.25
$ws.Range("P105").Value = -0.25
$ws.Range("Q105").Value = 1.86
$ws.Range("R105").Value = 2.04
$ws.Range("S105").Value = 2.75
$ws.Range("T105").Value = 1.975
$ws.Range("U105").Value = 1.875
$ws.Range("V105").Value = 1.15
$ws.Range("X105").Value = -1
$ws.Range("Y105").Value = 0.8600000000000001
$ws.Range("Z105").Value = -1
$ws.Range("AA105").Value = -1
$ws.Range("AB105").Value = 0.875

# Row 159 (id=157): becomes what row 160 used to contain
$ws.Range("B159").Value = 7127419
$ws.Range("E159").Value = "Wellington Phoenix"
$ws.Range("F159").Value = "Macarthur FC"
$ws.Range("G159").Value = 3
$ws.Range("H159").Value = 0
$ws.Range("I159").Value = "H"
$ws.Range("J159").Value = 1.85
$ws.Range("K159").Value = 3.5
$ws.Range("L159").Value = 3.9
$ws.Range("M159").Value = 1.55
$ws.Range("N159").Value = 4.5
$ws.Range("O159").Value = 5.25
$ws.Range("P159").Value = -1
$ws.Range("Q159").Value = 1.89
$ws.Range("R159").Value = 2.01
$ws.Range("S159").Value = 3.5
$ws.Range("T159").Value = 1.9
$ws.Range("U159").Value = 1.95
$ws.Range("V159").Value = 0.55
$ws.Range("X159").Value = -1
$ws.Range("Y159").Value = 0.8899999999999999
$ws.Range("Z159").Value = -1
$ws.Range("AA159").Value = -1
$ws.Range("AB159").Value = 0.95

# Row 160 (id=158): becomes what row 159 used to contain
$ws.Range("B160").Value = 7127418
$ws.Range("E160").Value = "Newcastle Jets"
$ws.Range("F160").Value = "Central Coast Mariners"
$ws.Range("G160").Value = 1
$ws.Range("H160").Value = 3
$ws.Range("I160").Value = "A"
$ws.Range("J160").Value = 3.6
$ws.Range("K160").Value = 3.25
$ws.Range("L160").Value = 2
$ws.Range("M160").Value = 4.2
$ws.Range("N160").Value = 4
$ws.Range("O160").Value = 1.75
$ws.Range("P160").Value = 0.75
$ws.Range("Q160").Value = 1.85
$ws.Range("R160").Value = 2
$ws.Range("S160").Value = 3
$ws.Range("T160").Value = 1.975
$ws.Range("U160").Value = 1.875
$ws.Range("V160").Value = -1
$ws.Range("X160").Value = 0.75
$ws.Range("Y160").Value = -1
$ws.Range("Z160").Value = 1
$ws.Range("AA160").Value = 0.9750000000000001
$ws.Range("AB160").Value = -1

# Row 171 (id=169): odds update
$ws.Range("M171").Value = 2.2
$ws.Range("Q171").Value = 1.91
$ws.Range("R171").Value = 1.99
